$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.471.97"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "1.900.71"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4907"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2923"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06677"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "1.894.80"
$ws.Range("E10").Value = "  +1.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6680"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.27%  "

$ws.Range("D16").Value = "30.454.84"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007884"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.419"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.65%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.130.44"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "194.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.139"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.498"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.942"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.488"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.334"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09158"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05148"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7387"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.728"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01843"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9258"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4413"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.902"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.78%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1372"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.586"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.079"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05836"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3938"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.87%  "
